# Appointment_Details.xlsx - add two new call-queue appointment rows (3 & 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Addr, [string]$Val)
    # Force the cell to stay TEXT even when the value looks numeric
    # (mirrors Excel's own quote-prefix / "Text" behaviour).
    if ($Val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Sheet.Range($Addr).Value = "'" + $Val
    } else {
        $Sheet.Range($Addr).Value = $Val
    }
}

# Row 3
Set-TextCell $ws "A3" "Aarav Mehta"
Set-TextCell $ws "B3" "Monday"
Set-TextCell $ws "C3" "morning"
Set-TextCell $ws "D3" "डॉ. से"
Set-TextCell $ws "E3" "28"
Set-TextCell $ws "F3" "Male"
Set-TextCell $ws "G3" "917823844614"
Set-TextCell $ws "H3" "24 MG Road, Bengaluru"
Set-TextCell $ws "I3" "2025-06-25 18:43:54"

# Row 4
Set-TextCell $ws "A4" "Aarav Mehta"
Set-TextCell $ws "B4" "Monday"
Set-TextCell $ws "C4" "morning"
Set-TextCell $ws "D4" "डॉ. से"
Set-TextCell $ws "E4" "28"
Set-TextCell $ws "F4" "Male"
Set-TextCell $ws "G4" "917823844614"
Set-TextCell $ws "H4" "24 MG Road, Bengaluru"
Set-TextCell $ws "I4" "2025-06-25 19:32:35"
